$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.11260631656082
$ws.Range("D2").Value = 7.701134910730441
$ws.Range("E2").Value = 14.96722265627409
$ws.Range("F2").Value = 45.68285142200519
$ws.Range("G2").Value = 55.45190140825297
$ws.Range("H2").Value = 20.91589483817961
$ws.Range("J2").Value = 11.89865745306942
$ws.Range("K2").Value = 10.59966846107984
$ws.Range("L2").Value = 9.277229157312012
$ws.Range("M2").Value = 16.07145015089271
$ws.Range("B3").Value = 17.06266031677161
$ws.Range("D3").Value = 7.708173286225478
$ws.Range("E3").Value = 14.98326580100813
$ws.Range("F3").Value = 45.6860255700572
$ws.Range("G3").Value = 55.36769166829395
$ws.Range("H3").Value = 20.94120542516939
$ws.Range("J3").Value = 11.90977629097207
$ws.Range("K3").Value = 10.45216576039691
$ws.Range("L3").Value = 9.275998421076055
$ws.Range("M3").Value = 16.07332847837192
$ws.Range("B4").Value = 17.03534555965846
$ws.Range("D4").Value = 7.712979496534544
$ws.Range("E4").Value = 14.99375550644383
$ws.Range("F4").Value = 45.69606354608946
$ws.Range("G4").Value = 55.32752680191187
$ws.Range("H4").Value = 20.95976722412843
$ws.Range("J4").Value = 11.91696071866224
$ws.Range("K4").Value = 10.36339746320939
$ws.Range("L4").Value = 9.276698121696993
$ws.Range("M4").Value = 16.0769411896587
$ws.Range("B5").Value = 17.02506645370864
$ws.Range("D5").Value = 7.715060247211788
$ws.Range("E5").Value = 14.99819121779132
$ws.Range("F5").Value = 45.70218850668409
$ws.Range("G5").Value = 55.31406864774587
$ws.Range("H5").Value = 20.96809054458543
$ws.Range("J5").Value = 11.91997854209577
$ws.Range("K5").Value = 10.32772417912428
$ws.Range("L5").Value = 9.277350421523362
$ws.Range("M5").Value = 16.07903352165654
$ws.Range("B6").Value = 17.02341131162967
$ws.Range("D6").Value = 7.715413143083749
$ws.Range("E6").Value = 14.99893750296236
$ws.Range("F6").Value = 45.70332841816473
$ws.Range("G6").Value = 55.31200982203225
$ws.Range("H6").Value = 20.96951846701046
$ws.Range("J6").Value = 11.92048509878519
$ws.Range("K6").Value = 10.32183227874754
$ws.Range("L6").Value = 9.277480943745299
$ws.Range("M6").Value = 16.07941845095736
$ws.Range("B7").Value = 17.03520347165838
$ws.Range("D7").Value = 7.713007063122927
$ws.Range("E7").Value = 14.99381467534785
$ws.Range("F7").Value = 45.6961379124683
$ws.Range("G7").Value = 55.32733351262852
$ws.Range("H7").Value = 20.95987640189841
$ws.Range("J7").Value = 11.91700105291406
$ws.Range("K7").Value = 10.36291427122951
$ws.Range("L7").Value = 9.276705430717504
$ws.Range("M7").Value = 16.07696689465525
$ws.Range("B8").Value = 17.09469476626296
$ws.Range("D8").Value = 7.703461351212066
$ws.Range("E8").Value = 14.97262188870916
$ws.Range("F8").Value = 45.6822674369306
$ws.Range("G8").Value = 55.42047586114912
$ws.Range("H8").Value = 20.9239947024512
$ws.Range("J8").Value = 11.90241718892398
$ws.Range("K8").Value = 10.54846526769926
$ws.Range("L8").Value = 9.276503493439378
$ws.Range("M8").Value = 16.07158841370327
$ws.Range("B9").Value = 17.23752789928453
$ws.Range("D9").Value = 7.688573629346819
$ws.Range("E9").Value = 14.93611841722837
$ws.Range("F9").Value = 45.71920561747216
$ws.Range("G9").Value = 55.69427166652435
$ws.Range("H9").Value = 20.87761676898121
$ws.Range("J9").Value = 11.87664336562698
$ws.Range("K9").Value = 10.92445973886255
$ws.Range("L9").Value = 9.287597922465666
$ws.Range("M9").Value = 16.08047496289758
$ws.Range("B10").Value = 17.35779754607266
$ws.Range("D10").Value = 7.679953001910225
$ws.Range("E10").Value = 14.9123594127291
$ws.Range("F10").Value = 45.78534073232168
$ws.Range("G10").Value = 55.95028892134268
$ws.Range("H10").Value = 20.85818245405025
$ws.Range("J10").Value = 11.85941426570864
$ws.Range("K10").Value = 11.20521616447804
$ws.Range("L10").Value = 9.302670096569814
$ws.Range("M10").Value = 16.09873330339255
$ws.Range("B11").Value = 17.41569127371415
$ws.Range("D11").Value = 7.676530580398931
$ws.Range("E11").Value = 14.90221068206482
$ws.Range("F11").Value = 45.82385458144198
$ws.Range("G11").Value = 56.07847174789789
$ws.Range("H11").Value = 20.85252117510425
$ws.Range("J11").Value = 11.85194368390005
$ws.Range("K11").Value = 11.33331617839133
$ws.Range("L11").Value = 9.311008978505324
$ws.Range("M11").Value = 16.10955704910727
$ws.Range("B12").Value = 17.43805806178887
$ws.Range("D12").Value = 7.67530605303471
$ws.Range("E12").Value = 14.89846208701609
$ws.Range("F12").Value = 45.83964561780225
$ws.Range("G12").Value = 56.12867576237252
$ws.Range("H12").Value = 20.85083438943647
$ws.Range("J12").Value = 11.84916730384843
$ws.Range("K12").Value = 11.38183023617407
$ws.Range("L12").Value = 9.314377863457615
$ws.Range("M12").Value = 16.11401491107441
$ws.Range("B13").Value = 17.43322146947439
$ws.Range("D13").Value = 7.675566603188272
$ws.Range("E13").Value = 14.89926521555959
$ws.Range("F13").Value = 45.83619116366992
$ws.Range("G13").Value = 56.11778978705839
$ws.Range("H13").Value = 20.85117734754808
$ws.Range("J13").Value = 11.84976291235697
$ws.Range("K13").Value = 11.37138240281015
$ws.Range("L13").Value = 9.313642957800127
$ws.Range("M13").Value = 16.11303890671126
$ws.Range("B14").Value = 17.41752259280307
$ws.Range("D14").Value = 7.676428407201784
$ws.Range("E14").Value = 14.90190039053383
$ws.Range("F14").Value = 45.82512956995573
$ws.Range("G14").Value = 56.08256884396622
$ws.Range("H14").Value = 20.8523732445426
$ws.Range("J14").Value = 11.85171421711515
$ws.Range("K14").Value = 11.33730763103959
$ws.Range("L14").Value = 9.311281919545445
$ws.Range("M14").Value = 16.10991661987877
$ws.Range("B15").Value = 17.40796393409587
$ws.Range("D15").Value = 7.676965585352447
$ws.Range("E15").Value = 14.9035268101302
$ws.Range("F15").Value = 45.81851099424364
$ws.Range("G15").Value = 56.0612110264045
$ws.Range("H15").Value = 20.85316527497184
$ws.Range("J15").Value = 11.8529162872126
$ws.Range("K15").Value = 11.31643509809216
$ws.Range("L15").Value = 9.309863147775991
$ws.Range("M15").Value = 16.10805080948174
$ws.Range("B16").Value = 17.35407691715429
$ws.Range("D16").Value = 7.680186683037396
$ws.Range("E16").Value = 14.91303589790335
$ws.Range("F16").Value = 45.78299293208254
$ws.Range("G16").Value = 55.94214592890337
$ws.Range("H16").Value = 20.85861639872765
$ws.Range("J16").Value = 11.85990985320997
$ws.Range("K16").Value = 11.19684795018354
$ws.Range("L16").Value = 9.302154790537154
$ws.Range("M16").Value = 16.0980763678684
$ws.Range("B17").Value = 17.32182437256762
$ws.Range("D17").Value = 7.682290346736083
$ws.Range("E17").Value = 14.91903807417978
$ws.Range("F17").Value = 45.76335902774686
$ws.Range("G17").Value = 55.87209053901832
$ws.Range("H17").Value = 20.86277474398247
$ws.Range("J17").Value = 11.86429403241752
$ws.Range("K17").Value = 11.12354734174245
$ws.Range("L17").Value = 9.29780424996574
$ws.Range("M17").Value = 16.09260025869755
$ws.Range("B18").Value = 17.303573988349
$ws.Range("D18").Value = 7.683547318775447
$ws.Range("E18").Value = 14.92255244739581
$ws.Range("F18").Value = 45.75285985097142
$ws.Range("G18").Value = 55.832901025714
$ws.Range("H18").Value = 20.86546580346244
$ws.Range("J18").Value = 11.86685025464889
$ws.Range("K18").Value = 11.08142566836072
$ws.Range("L18").Value = 9.295441603245912
$ws.Range("M18").Value = 16.08968773489125
$ws.Range("B19").Value = 17.29744673438179
$ws.Range("D19").Value = 7.683980989601122
$ws.Range("E19").Value = 14.92375302548368
$ws.Range("F19").Value = 45.74944148074363
$ws.Range("G19").Value = 55.81982243832508
$ws.Range("H19").Value = 20.86642835648008
$ws.Range("J19").Value = 11.86772168943317
$ws.Range("K19").Value = 11.06717220855812
$ws.Range("L19").Value = 9.294665699973777
$ws.Range("M19").Value = 16.08874242759931
$ws.Range("B20").Value = 17.32522671691224
$ws.Range("D20").Value = 7.682061546023461
$ws.Range("E20").Value = 14.91839270944136
$ws.Range("F20").Value = 45.76536697842138
$ws.Range("G20").Value = 55.87943389058007
$ws.Range("H20").Value = 20.86230110656243
$ws.Range("J20").Value = 11.86382375393131
$ws.Range("K20").Value = 11.13134665474796
$ws.Range("L20").Value = 9.298252931569035
$ws.Range("M20").Value = 16.09315867219259
$ws.Range("B21").Value = 17.42212180615933
$ws.Range("D21").Value = 7.676173337262998
$ws.Range("E21").Value = 14.90112381324293
$ws.Range("F21").Value = 45.82834592668281
$ws.Range("G21").Value = 56.09286910996545
$ws.Range("H21").Value = 20.85200957953606
$ws.Range("J21").Value = 11.85113964661571
$ws.Range("K21").Value = 11.34731646066719
$ws.Range("L21").Value = 9.311969700692913
$ws.Range("M21").Value = 16.11082398858387
$ws.Range("B22").Value = 17.48802759728706
$ws.Range("D22").Value = 7.672741488090401
$ws.Range("E22").Value = 14.8903883224113
$ws.Range("F22").Value = 45.87653659263502
$ws.Range("G22").Value = 56.242048753546
$ws.Range("H22").Value = 20.84794712502622
$ws.Range("J22").Value = 11.84315615467109
$ws.Range("K22").Value = 11.48847448088951
$ws.Range("L22").Value = 9.322163986859969
$ws.Range("M22").Value = 16.12446114778307
$ws.Range("B23").Value = 17.45262120666658
$ws.Range("D23").Value = 7.674535128182612
$ws.Range("E23").Value = 14.89606776422311
$ws.Range("F23").Value = 45.85017507144531
$ws.Range("G23").Value = 56.16154997804992
$ws.Range("H23").Value = 20.84987171109433
$ws.Range("J23").Value = 11.84738913653425
$ws.Range("K23").Value = 11.41315096930457
$ws.Range("L23").Value = 9.316611310594372
$ws.Range("M23").Value = 16.11699235966596
$ws.Range("B24").Value = 17.32368760632144
$ws.Range("D24").Value = 7.68216483871876
$ws.Range("E24").Value = 14.91868428059504
$ws.Range("F24").Value = 45.76445672694556
$ws.Range("G24").Value = 55.87611057599833
$ws.Range("H24").Value = 20.86251430233074
$ws.Range("J24").Value = 11.86403625562167
$ws.Range("K24").Value = 11.12782051977903
$ws.Range("L24").Value = 9.298049650928119
$ws.Range("M24").Value = 16.09290547884032
$ws.Range("B25").Value = 17.19614731403089
$ws.Range("D25").Value = 7.692192870319964
$ws.Range("E25").Value = 14.94545461993382
$ws.Range("F25").Value = 45.70235348236928
$ws.Range("G25").Value = 55.61050838294113
$ws.Range("H25").Value = 20.88759305429181
$ws.Range("J25").Value = 11.88331501356916
$ws.Range("K25").Value = 10.82174321009839
$ws.Range("L25").Value = 9.283374534288132
$ws.Range("M25").Value = 16.0760017191003
